# Docx writer: Use a different style for block quotes in notes.
#
# Adds a new paragraph style "Footnote Block Text" (styleId
# "FootnoteBlockText"), based on "Footnote Text" / followed-by
# "Footnote Text", mirroring the existing "Block Text" style but scoped
# to footnotes so it can later be given its own font size.

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$style = $d.Styles.Add("Footnote Block Text", 1)

$style.BaseStyle = "Footnote Text"
$style.NextParagraphStyle = "Footnote Text"

$style.Priority = 9
$style.UnhideWhenUsed = $true
$style.QuickStyle = $true

$style.ParagraphFormat.SpaceBefore = 5
$style.ParagraphFormat.SpaceAfter = 5
$style.ParagraphFormat.FirstLineIndent = 0
$style.ParagraphFormat.LeftIndent = 24
$style.ParagraphFormat.RightIndent = 24
